$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 345 (they become the new rows 345-346,
# and the old rows 345-351 shift down to become 347-353).
$ws.Rows.Item(345).Resize(2).Insert()

# Give column D in the two new rows the same date number format used by the
# other "Fecha" cells in this column (style index 2 / numFmtId 165).
$ws.Range("D345:D346").NumberFormat = $ws.Range("D347").NumberFormat

# --- New row 345 ---
$ws.Cells.Item(345, 1).Value = 1
$ws.Cells.Item(345, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(345, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(345, 4).Value = 44890
$ws.Cells.Item(345, 5).Value = 15
$ws.Cells.Item(345, 6).Value = 100114013
$ws.Cells.Item(345, 7).Value = "Zanahoria"
$ws.Cells.Item(345, 8).Value = "Sin especificar"
$ws.Cells.Item(345, 9).Value = "Primera"
$ws.Cells.Item(345, 10).Value = 150
$ws.Cells.Item(345, 11).Value = 32000
$ws.Cells.Item(345, 12).Value = 33000
$ws.Cells.Item(345, 13).Value = 32400
$ws.Cells.Item(345, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(345, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(345, 16).Value = 1296
$ws.Cells.Item(345, 17).Value = 25
$ws.Cells.Item(345, 18).Value = "Hortaliza"

# --- New row 346 ---
$ws.Cells.Item(346, 1).Value = 1
$ws.Cells.Item(346, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(346, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(346, 4).Value = 44890
$ws.Cells.Item(346, 5).Value = 15
$ws.Cells.Item(346, 6).Value = 100114013
$ws.Cells.Item(346, 7).Value = "Zanahoria"
$ws.Cells.Item(346, 8).Value = "Sin especificar"
$ws.Cells.Item(346, 9).Value = "Primera"
$ws.Cells.Item(346, 10).Value = 180
$ws.Cells.Item(346, 11).Value = 32000
$ws.Cells.Item(346, 12).Value = 33000
$ws.Cells.Item(346, 13).Value = 32333
$ws.Cells.Item(346, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(346, 15).Value = "Región de Tarapacá"
$ws.Cells.Item(346, 16).Value = 1293
$ws.Cells.Item(346, 17).Value = 25
$ws.Cells.Item(346, 18).Value = "Hortaliza"
